$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Record_PlayerHero")

# --- Row 2 data change: Col value 21 -> 31 -----------------------------
$ws.Range("C2").Value2 = 31

# --- New header row (row 1): Skill1..Skill5, SkillLevel1..SkillLevel5 --
$ws.Range("AG1").Value2 = "Skill1"
$ws.Range("AH1").Value2 = "Skill2"
$ws.Range("AI1").Value2 = "Skill3"
$ws.Range("AJ1").Value2 = "Skill4"
$ws.Range("AK1").Value2 = "Skill5"
$ws.Range("AL1").Value2 = "SkillLevel1"
$ws.Range("AM1").Value2 = "SkillLevel2"
$ws.Range("AN1").Value2 = "SkillLevel3"
$ws.Range("AO1").Value2 = "SkillLevel4"
$ws.Range("AP1").Value2 = "SkillLevel5"

# match the formatting already used by the other header cells (W1:AF1)
$ws.Range("AG1:AP1").Font.Size = $ws.Range("X1").Font.Size

# --- New type row (row 2): string x5, int x5 ----------------------------
$ws.Range("AG2").Value2 = "string"
$ws.Range("AH2").Value2 = "string"
$ws.Range("AI2").Value2 = "string"
$ws.Range("AJ2").Value2 = "string"
$ws.Range("AK2").Value2 = "string"
$ws.Range("AL2").Value2 = "int"
$ws.Range("AM2").Value2 = "int"
$ws.Range("AN2").Value2 = "int"
$ws.Range("AO2").Value2 = "int"
$ws.Range("AP2").Value2 = "int"

# --- New column width for SkillLevel1 (col 38 / AL), matching bestFit --
$ws.Columns.Item(38).ColumnWidth = 11.625

# --- sheet view: drop the frozen/scrolled "K1" top-left cell, select B2
$ws.Activate()
$excel.ActiveWindow.ScrollColumn = 1
$excel.ActiveWindow.ScrollRow = 1
$ws.Range("B2").Select()
